$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Peru Liga 1")

# Cell-level updates: the underlying "match" rows for several fixtures were
# re-paired with the correct external match id (column B) during the base
# refresh, which re-shuffles every odds column (F..AC) between the affected
# rows. Column A (running index) and E (kickoff date, shared between the
# swapped fixtures) are left untouched.
$updates = @(
    @(156, 2, 7211640),
    @(156, 6, "UTC Cajamarca"),
    @(156, 7, "Sport Boys"),
    @(156, 8, 1),
    @(156, 9, 1),
    @(156, 10, "D"),
    @(156, 11, 1.615),
    @(156, 12, 3.75),
    @(156, 13, 5),
    @(156, 14, 1.5),
    @(156, 15, 4.2),
    @(156, 16, 6.5),
    @(156, 17, -1),
    @(156, 18, 1.8),
    @(156, 19, 2.05),
    @(156, 20, 2.5),
    @(156, 21, 1.875),
    @(156, 22, 1.975),
    @(156, 23, -1),
    @(156, 24, 3.2),
    @(156, 26, -1),
    @(156, 27, 1.05),
    @(156, 29, 0.9750000000000001),
    @(157, 2, 7211641),
    @(157, 6, "Sport Huancayo"),
    @(157, 7, "Deportivo Municipal"),
    @(157, 8, 2),
    @(157, 9, 0),
    @(157, 10, "H"),
    @(157, 11, 1.125),
    @(157, 12, 7),
    @(157, 13, 17),
    @(157, 14, 1.166),
    @(157, 15, 6.5),
    @(157, 16, 12),
    @(157, 17, -2),
    @(157, 18, 1.775),
    @(157, 19, 2.025),
    @(157, 20, 3.5),
    @(157, 21, 1.9),
    @(157, 22, 1.9),
    @(157, 23, 0.1659999999999999),
    @(157, 24, -1),
    @(157, 26, 0),
    @(157, 27, -0),
    @(157, 29, 0.8999999999999999),
    @(175, 2, 7302795),
    @(175, 6, "Unin Comercio"),
    @(175, 7, "Deportivo Garcilaso"),
    @(175, 9, 2),
    @(175, 10, "A"),
    @(175, 11, 2.25),
    @(175, 12, 3.3),
    @(175, 13, 2.7),
    @(175, 14, 1.75),
    @(175, 15, 3.6),
    @(175, 16, 4),
    @(175, 17, -0.5),
    @(175, 18, 1.8),
    @(175, 19, 2),
    @(175, 20, 2.75),
    @(175, 21, 1.825),
    @(175, 22, 1.975),
    @(175, 23, -1),
    @(175, 25, 3),
    @(175, 27, 1),
    @(175, 28, 0.4125),
    @(175, 29, -0.5),
    @(176, 2, 7302796),
    @(176, 6, "Sport Huancayo"),
    @(176, 7, "Sport Boys"),
    @(176, 9, 0),
    @(176, 10, "H"),
    @(176, 11, 1.727),
    @(176, 12, 3.75),
    @(176, 13, 4.333),
    @(176, 14, 1.25),
    @(176, 15, 5.25),
    @(176, 16, 10),
    @(176, 17, -1.75),
    @(176, 18, 1.925),
    @(176, 19, 1.875),
    @(176, 20, 3),
    @(176, 21, 1.875),
    @(176, 22, 1.925),
    @(176, 23, 0.25),
    @(176, 25, -1),
    @(176, 27, 0.875),
    @(176, 28, -1),
    @(176, 29, 0.925),
    @(180, 2, 7384622),
    @(180, 6, "Deportivo Municipal"),
    @(180, 7, "Academia Deportiva Cantolao"),
    @(180, 8, 1),
    @(180, 9, 2),
    @(180, 10, "A"),
    @(180, 11, 1.444),
    @(180, 12, 4.333),
    @(180, 13, 7),
    @(180, 14, 1.5),
    @(180, 15, 3.75),
    @(180, 16, 6),
    @(180, 18, 1.825),
    @(180, 19, 2.025),
    @(180, 20, 2.75),
    @(180, 21, 1.875),
    @(180, 22, 1.975),
    @(180, 23, -1),
    @(180, 25, 5),
    @(180, 26, -1),
    @(180, 27, 1.025),
    @(180, 28, 0.4375),
    @(180, 29, -0.5),
    @(181, 2, 7384623),
    @(181, 6, "Sport Boys"),
    @(181, 7, "Cienciano"),
    @(181, 8, 2),
    @(181, 9, 1),
    @(181, 10, "H"),
    @(181, 11, 2),
    @(181, 12, 3.4),
    @(181, 13, 3.5),
    @(181, 14, 1.833),
    @(181, 15, 4),
    @(181, 16, 3.2),
    @(181, 17, -0.5),
    @(181, 18, 1.925),
    @(181, 19, 1.875),
    @(181, 20, 3),
    @(181, 21, 1.925),
    @(181, 22, 1.875),
    @(181, 23, 0.833),
    @(181, 25, -1),
    @(181, 26, 0.925),
    @(181, 27, -1),
    @(181, 28, 0),
    @(181, 29, -0),
    @(182, 2, 7384624),
    @(182, 6, "Cesar Vallejo"),
    @(182, 7, "Cusco FC"),
    @(182, 8, 3),
    @(182, 14, 1.45),
    @(182, 15, 4.2),
    @(182, 16, 6.5),
    @(182, 17, -1),
    @(182, 18, 1.75),
    @(182, 19, 2.05),
    @(182, 20, 2.5),
    @(182, 21, 1.95),
    @(182, 22, 1.85),
    @(182, 23, 0.45),
    @(182, 26, 0.75),
    @(182, 28, 0.95),
    @(182, 29, -1),
    @(183, 2, 7384630),
    @(183, 6, "Atletico Grau"),
    @(183, 7, "Unin Comercio"),
    @(183, 8, 0),
    @(183, 9, 1),
    @(183, 11, 2.8),
    @(183, 12, 3.4),
    @(183, 13, 2.15),
    @(183, 14, 1.75),
    @(183, 16, 3.8),
    @(183, 17, -0.75),
    @(183, 18, 2),
    @(183, 19, 1.8),
    @(183, 20, 3),
    @(183, 21, 1.85),
    @(183, 22, 1.95),
    @(183, 25, 2.8),
    @(183, 27, 0.8),
    @(183, 28, -1),
    @(183, 29, 0.95),
    @(184, 2, 7384625),
    @(184, 6, "AD Tarma"),
    @(184, 7, "Carlos Manucci"),
    @(184, 9, 0),
    @(184, 10, "D"),
    @(184, 11, 1.5),
    @(184, 12, 3.75),
    @(184, 13, 7),
    @(184, 14, 1.363),
    @(184, 15, 4.333),
    @(184, 16, 9.5),
    @(184, 17, -1.25),
    @(184, 18, 1.875),
    @(184, 19, 1.925),
    @(184, 20, 2.5),
    @(184, 21, 1.8),
    @(184, 22, 2),
    @(184, 24, 3.333),
    @(184, 25, -1),
    @(184, 27, 0.925),
    @(184, 29, 1),
    @(185, 2, 7384628),
    @(185, 6, "Deportivo Binacional"),
    @(185, 7, "FBC Melgar"),
    @(185, 8, 1),
    @(185, 9, 2),
    @(185, 10, "A"),
    @(185, 11, 2.75),
    @(185, 12, 3.3),
    @(185, 13, 2.375),
    @(185, 14, 3.3),
    @(185, 15, 3.6),
    @(185, 16, 2),
    @(185, 17, 0.5),
    @(185, 18, 1.8),
    @(185, 19, 2),
    @(185, 20, 2.75),
    @(185, 21, 1.975),
    @(185, 22, 1.875),
    @(185, 24, -1),
    @(185, 25, 1),
    @(185, 27, 1),
    @(185, 28, 0.4875),
    @(185, 29, -0.5),
    @(186, 2, 7384627),
    @(186, 6, "Universitario de Deportes"),
    @(186, 7, "Sport Huancayo"),
    @(186, 8, 2),
    @(186, 11, 1.25),
    @(186, 13, 12),
    @(186, 14, 1.181),
    @(186, 15, 6),
    @(186, 17, -1.75),
    @(186, 18, 1.8),
    @(186, 19, 2),
    @(186, 20, 2.75),
    @(186, 21, 1.85),
    @(186, 22, 1.95),
    @(186, 23, 0.181),
    @(186, 26, 0.4),
    @(186, 27, -0.5),
    @(186, 28, -1),
    @(186, 29, 0.95),
    @(188, 2, 7384626),
    @(188, 6, "Sporting Cristal"),
    @(188, 7, "Alianza Atletico"),
    @(188, 8, 3),
    @(188, 11, 1.3),
    @(188, 13, 9),
    @(188, 14, 1.166),
    @(188, 15, 6.5),
    @(188, 17, -2),
    @(188, 18, 1.85),
    @(188, 19, 1.95),
    @(188, 20, 3.25),
    @(188, 21, 2),
    @(188, 22, 1.8),
    @(188, 23, 0.1659999999999999),
    @(188, 26, 0.8500000000000001),
    @(188, 27, -1),
    @(188, 28, -0.5),
    @(188, 29, 0.4),
    @(228, 2, 7818817),
    @(228, 6, "Sport Boys"),
    @(228, 7, "Cusco FC"),
    @(228, 8, 3),
    @(228, 10, "H"),
    @(228, 11, 2.2),
    @(228, 12, 3.2),
    @(228, 13, 3.2),
    @(228, 14, 1.6),
    @(228, 15, 3.75),
    @(228, 16, 5.75),
    @(228, 17, -0.75),
    @(228, 18, 1.85),
    @(228, 19, 2),
    @(228, 20, 2.5),
    @(228, 21, 1.975),
    @(228, 22, 1.875),
    @(228, 23, 0.6000000000000001),
    @(228, 24, -1),
    @(228, 26, 0.8500000000000001),
    @(228, 28, 0.9750000000000001),
    @(228, 29, -1),
    @(229, 2, 7818816),
    @(229, 6, "UTC Cajamarca"),
    @(229, 7, "Universitario de Deportes"),
    @(229, 8, 0),
    @(229, 10, "D"),
    @(229, 11, 3.3),
    @(229, 12, 3.3),
    @(229, 13, 2.1),
    @(229, 14, 4.5),
    @(229, 15, 3.2),
    @(229, 16, 1.95),
    @(229, 17, 0.5),
    @(229, 18, 2),
    @(229, 19, 1.85),
    @(229, 20, 2),
    @(229, 21, 1.775),
    @(229, 22, 2.1),
    @(229, 23, -1),
    @(229, 24, 2.2),
    @(229, 26, 1),
    @(229, 28, -1),
    @(229, 29, 1.1)
)

foreach ($u in $updates) {
    $r = $u[0]
    $c = $u[1]
    $v = $u[2]
    $ws.Cells.Item($r, $c).Value = $v
}

# New fixture appended at the bottom of the table (row 281).
$lastRow = 280
$newRow = 281

$ws.Cells.Item($lastRow, 1).Copy() | Out-Null
$ws.Cells.Item($newRow, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item($lastRow, 5).Copy() | Out-Null
$ws.Cells.Item($newRow, 5).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Cells.Item($newRow, 1).Value = 279
$ws.Cells.Item($newRow, 2).Value = 8042070
$ws.Cells.Item($newRow, 3).Value = "Peru Liga 1"
$ws.Cells.Item($newRow, 4).Value = "Peru Liga 1"
$ws.Cells.Item($newRow, 5).Value = 45394.70833333334
$ws.Cells.Item($newRow, 6).Value = "Sport Huancayo"
$ws.Cells.Item($newRow, 7).Value = "Cienciano"
$ws.Cells.Item($newRow, 11).Value = 1.8
$ws.Cells.Item($newRow, 12).Value = 3.5
$ws.Cells.Item($newRow, 13).Value = 4.333
$ws.Cells.Item($newRow, 14).Value = 1.95
$ws.Cells.Item($newRow, 15).Value = 3.4
$ws.Cells.Item($newRow, 16).Value = 3.8
$ws.Cells.Item($newRow, 17).Value = -0.5
$ws.Cells.Item($newRow, 18).Value = 2.025
$ws.Cells.Item($newRow, 19).Value = 1.825
$ws.Cells.Item($newRow, 20).Value = 2.5
$ws.Cells.Item($newRow, 21).Value = 1.875
$ws.Cells.Item($newRow, 22).Value = 1.975
$ws.Cells.Item($newRow, 23).Value = 0
$ws.Cells.Item($newRow, 24).Value = 0
$ws.Cells.Item($newRow, 25).Value = 0
$ws.Cells.Item($newRow, 26).Value = 0
$ws.Cells.Item($newRow, 27).Value = 0
